# Menu.xlsx — "Adds ID and store lookup methods to StoreData" / fixes a
# missing param on the drinks menu sheet: the Amount_2 (column G) option
# list for "Ah Beng Drink" needs a 4th value to line up with the
# Normal,Chino,O,Gao option list in column F.

$wb = $excel.ActiveWorkbook

$drinkSheet = $wb.Worksheets.Item("Ah Beng Drink")
$foodSheet  = $wb.Worksheets.Item("Ah Lian food")

# Fix the missing param: column G ("Amount_2") only had 3 comma-separated
# amounts but needs 4 (to match Option_2's "Normal,Chino,O,Gao").
$drinkSheet.Range("G2").Value = "0.00,0.00,0.00,0.00"
$drinkSheet.Range("G3").Value = "0.00,0.00,0.00,0.00"
$drinkSheet.Range("G4").Value = "0.00,0.00,0.00,0.00"

# Make the drinks sheet the active tab, with G2:G4 selected (the range
# that was just fixed) and G2 as the active cell.
$drinkSheet.Activate()
$drinkSheet.Range("G2:G4").Select()
